$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 61, shifting the existing rows 61..161 down to 62..162
$ws.Rows(61).Insert()

# Populate the newly inserted row 61 with the new weekly record
$ws.Range("A61").Value = 11
$ws.Range("B61").Value = "Vega Monumental Concepción"
$ws.Range("C61").Value = "Bíobío"
$ws.Range("D61").Value = 44495
$ws.Range("E61").Value = 8
$ws.Range("F61").Value = 100112008
$ws.Range("G61").Value = "Coliflor"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 2400
$ws.Range("K61").Value = 700
$ws.Range("L61").Value = 750
$ws.Range("M61").Value = 725
$ws.Range("N61").Value = "$/unidad"
$ws.Range("O61").Value = "Región Metropolitana"
$ws.Range("P61").Value = 725
$ws.Range("Q61").Value = 1
$ws.Range("R61").Value = "Hortaliza"
